$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 10 (only row 1 should remain)
$ws.Range("A2:B10").EntireRow.Delete() | Out-Null

# Update remaining row 1 values
$ws.Range("A1").Value = "Manchester City v Aston Villa"
$ws.Range("B1").Value = "had Stadium"
